$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns:
#  - one before the old "Start Date" column (F) to hold "Percentage Previous Year"
#  - one before the old "Value" -> "Year" boundary (after the shift, at J) to hold "Percentage"
$ws.Columns("F").Insert()
$ws.Columns("J").Insert()

# New headers
$ws.Range("F1").Value = "Percentage Previous Year"
$ws.Range("J1").Value = "Percentage"

# New data values (manual turnover-rent percentages)
$ws.Range("F2").Value = 6
$ws.Range("J2").Value = 6.25

$ws.Range("F3").Value = 1.25
$ws.Range("J3").Value = 1.3

# Re-apply column widths (existing columns shifted / resized, two new narrow columns added)
$ws.Columns("B").ColumnWidth = 18.666666666666668
$ws.Columns("C").ColumnWidth = 18.333333333333332
$ws.Columns("D").ColumnWidth = 17.5
$ws.Columns("E").ColumnWidth = 14.833333333333334
$ws.Columns("F").ColumnWidth = 19.166666666666668
$ws.Columns("G").ColumnWidth = 9.666666666666666
$ws.Columns("H").ColumnWidth = 9.333333333333334
$ws.Columns("I").ColumnWidth = 5.333333333333333
$ws.Columns("J").ColumnWidth = 5.333333333333333
$ws.Columns("K").ColumnWidth = 4.333333333333333

# Final selection left by the editor
$ws.Range("J4").Select() | Out-Null
